# cronograma_curso_provas.xlsx - "teste unidade 1 - modelos informacionais"
#
# Record the scores for the "Modelos Informacionais" course (row 4) under the
# "Unidade 1" activity columns (F:I). The weighted-total formula in column K
# and the class-average formula in K16 already exist on the sheet and will
# recalculate automatically once the inputs are present.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

$ws.Range("F4:I4").Value = 10

# Match the workbook's saved cursor position/selection (the merged activity
# row for "Unidade 1" lower on the sheet).
$ws.Range("E18:K18").Select()
